$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.4
$ws.Range("L2").Value = 3
$ws.Range("P2").Value = 3.95
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 2.32
$ws.Range("W2").Value = 12
$ws.Range("AC2").Value = 15
$ws.Range("AF2").Value = 34
$ws.Range("AK2").Value = 23
$ws.Range("AN2").Value = 5
$ws.Range("AT2").Value = 3.2
$ws.Range("AU2").Value = 7
$ws.Range("AY2").Value = 19
$ws.Range("A3").Value = "lleFOanB"
$ws.Range("C3").Value = "11:45"
$ws.Range("D3").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E3").Value = "Al Feiha"
$ws.Range("F3").Value = "Al Orubah"
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.25
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.75
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 10
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = 9.5
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 9.5
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 29
$ws.Range("AM3").Value = 34
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 12
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 51
$ws.Range("AT3").Value = 2.63
$ws.Range("AU3").Value = 8
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 19
$ws.Range("AY3").Value = 29
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 81
$ws.Range("BB3").Value = 300
$ws.Range("BC3").Value = 81
$ws.Range("BD3").Value = 81
$ws.Range("A4").Value = "vyTPDrYh"
$ws.Range("C4").Value = "12:05"
$ws.Range("E4").Value = "Al Okhdood"
$ws.Range("F4").Value = "Al Kholood"
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.92
$ws.Range("R4").Value = 1.82
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("X4").Value = 9.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 7
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 13
$ws.Range("AK4").Value = 41
$ws.Range("AO4").Value = 10
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 34
$ws.Range("AT4").Value = 2.75
$ws.Range("AX4").Value = 21
$ws.Range("AZ4").Value = 67
$ws.Range("BB4").Value = 400
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
